# Portchannels and logical device build
# Fills in the Port-Channel Configuration table (rows 56-59) on both the
# "FXOS DC1 Settings" and "FXOS DC2 Settings" worksheets, and adds the new
# Logical Device Configuration table (rows 61-62) on "FXOS DC2 Settings".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # FXOS DC1 Settings
$ws2 = $wb.Worksheets.Item(2)   # FXOS DC2 Settings

# ---------------------------------------------------------------------
# Helper data: the Port-Channel Configuration table is identical on both
# sheets except for the Port Type of Portchannel 30 (data vs mgmt).
# ---------------------------------------------------------------------

function Fill-PortchannelTable {
    param($ws, [string]$pc30PortType)

    # Row 56 - bold header row
    $ws.Range("B56:O56").Font.Bold = $true
    $ws.Range("B56").Value = "ID"
    $ws.Range("C56").Value = "Port Type"
    $ws.Range("D56").Value = "Auto-negotiation"
    $ws.Range("E56").Value = "Speed"
    $ws.Range("F56").Value = "Duplex"
    $ws.Range("G56").Value = "Port-channel Mode"
    $ws.Range("H56:O56").Value = "Member Interface"

    # Row 57 - Portchannel 10
    $ws.Range("A57").Value = "Portchannel"
    $ws.Range("B57").Value = 10
    $ws.Range("C57").Value = "data"
    $ws.Range("D57").Value = "off"
    $ws.Range("E57").Value = "1gbps"
    $ws.Range("F57").Value = "fullduplex"
    $ws.Range("G57").Value = "on"
    $ws.Range("H57").Value = "Ethernet1/1"
    $ws.Range("I57").Value = "Ethernet1/2"
    $ws.Range("J57").Value = "Ethernet1/3"
    $ws.Range("K57").Value = "Ethernet1/4"
    $ws.Range("L57").Value = "Ethernet1/5"
    $ws.Range("M57").Value = "Ethernet1/6"
    $ws.Range("N57").Value = "Ethernet1/7"
    $ws.Range("O57").Value = "Ethernet1/8"

    # Row 58 - Portchannel 20
    $ws.Range("A58").Value = "Portchannel"
    $ws.Range("B58").Value = 20
    $ws.Range("C58").Value = "data"
    $ws.Range("D58").Value = "off"
    $ws.Range("E58").Value = "1gbps"
    $ws.Range("F58").Value = "fullduplex"
    $ws.Range("G58").Value = "on"
    $ws.Range("H58").Value = "Ethernet2/1"
    $ws.Range("I58").Value = "Ethernet2/2"
    $ws.Range("J58").Value = "Ethernet2/3"
    $ws.Range("K58").Value = "Ethernet2/4"

    # Row 59 - Portchannel 30 (new row)
    $ws.Range("A59").Value = "Portchannel"
    $ws.Range("B59").Value = 30
    $ws.Range("C59").Value = $pc30PortType
    $ws.Range("D59").Value = "off"
    $ws.Range("E59").Value = "1gbps"
    $ws.Range("F59").Value = "fullduplex"
    $ws.Range("G59").Value = "on"
    $ws.Range("H59").Value = "Ethernet2/5"
    $ws.Range("I59").Value = "Ethernet2/6"
}

# FXOS DC1 Settings: Portchannel 30 is a data port
Fill-PortchannelTable $ws1 "data"

# FXOS DC2 Settings: Portchannel 30 is a mgmt port
Fill-PortchannelTable $ws2 "mgmt"

# ---------------------------------------------------------------------
# FXOS DC2 Settings only: new "# Logical Device Configuration" table
# ---------------------------------------------------------------------

$ws2.Range("A61:G61").Font.Bold = $true

$ws2.Range("A61").Value = "# Logical Device Configuration"
$ws2.Range("A62").Value = "ASA"
$ws2.Range("B61").Value = "Slot Number"
$ws2.Range("B62").Value = 1
$ws2.Range("C61").Value = "Hostname"
$ws2.Range("D61").Value = "Software Version"
$ws2.Range("C62").Value = "ASA1"
$ws2.Range("D62").Value = "9.12.1"
$ws2.Range("E61").Value = "Management Interface"
$ws2.Range("E62").Value = "Portchannel30"
$ws2.Range("F61").Value = "Nameif"
$ws2.Range("F62").Value = "management"
$ws2.Range("G61").Value = "Description"
$ws2.Range("G62").Value = "management link"

# ---------------------------------------------------------------------
# Column width tweaks on FXOS DC2 Settings (best-effort approximation;
# the COM layer quantizes widths, so these land close to, not exactly
# on, the original author's autofit values).
# ---------------------------------------------------------------------

$ws2.Columns.Item(7).ColumnWidth = 17.59
$ws2.Range("I1:J1").EntireColumn.ColumnWidth = 16.59
$ws2.Range("L1:N1").EntireColumn.ColumnWidth = 16.59

# ---------------------------------------------------------------------
# Selection / active-cell bookkeeping to mirror the final view state.
# FXOS DC1 Settings keeps its current scroll position; only the
# selection moves. FXOS DC2 Settings is reselected last so it remains
# the active tab.
# ---------------------------------------------------------------------

$ws1.Range("A70").Select()
$ws2.Range("G72").Select()
